$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Antwoordenblad")

$ws.Range("D16").Value = 3.100992589961451
$ws.Range("D37").Value = 70.91120726548422
$ws.Range("D38").Value = -0.0004327642678354389
$ws.Range("D39").Value = 7.104063194937617
$ws.Range("D49").Value = -2.980232238769531 / 100000000
$ws.Range("D50").Value = -6563531.254794854
$ws.Range("D51").Value = 104018.5367347449
$ws.Range("D57").Value = 32
$ws.Range("D58").Value = -2
$ws.Range("D59").Value = 25
$ws.Range("D69").Value = 4
$ws.Range("D77").Value = 44.146
$ws.Range("D81").Value = 16
$ws.Range("D84").Value = 16.13
$ws.Range("D85").Value = 30.073
$ws.Range("D89").Value = 44.146
